# Update "想去人数" (F column) figures and one cover image URL (I column)
# for the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, matching the
# regenerated site output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row number -> new F-column value
$updates = @{
    "展览" = @{
        5  = 1827
        9  = 2381
        10 = 133
        11 = 73
        13 = 1442
        16 = 316
        17 = 222
        19 = 177
        24 = 99
        25 = 35
        26 = 1488
        28 = 373
        29 = 261
        30 = 183
        31 = 287
        32 = 373
    }
    "全部类型" = @{
        5  = 1827
        10 = 2381
        11 = 133
        12 = 73
        14 = 1442
        17 = 316
        18 = 222
        20 = 177
        25 = 99
        26 = 35
        27 = 1488
        29 = 373
        30 = 261
        31 = 183
        32 = 287
        33 = 373
    }
}

# Rows whose Cover (I column) URL changed, per sheet.
$coverUpdates = @{
    "展览"   = @{ 29 = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png" }
    "全部类型" = @{ 30 = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates[$sheetName].Keys) {
        $newValue = $updates[$sheetName][$row]
        $ws.Range("F$row").Value = $newValue
    }

    if ($coverUpdates.ContainsKey($sheetName)) {
        foreach ($row in $coverUpdates[$sheetName].Keys) {
            $newUrl = $coverUpdates[$sheetName][$row]
            $ws.Range("I$row").Value = $newUrl
        }
    }
}
